# Update the "Förändrad" (Changed) date column for rows 2-6
# from 2023-10-09 (45208) to 2023-10-13 (45212)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = 45212
}
